# Daily scrape update - 2025-10-09 03:10:53 UTC
# Updates the Global Talent opportunities table: refreshes rows 2-13 with the
# latest scrape data, drops the now-stale last row (14), restores the E9 cell
# to the default (unstyled) format, and re-tunes a few column widths.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Delete the last data row (14) - it fell out of the latest scrape window.
#    This also shrinks the sheet dimension from A1:H14 to A1:H13.
# ---------------------------------------------------------------------------
$ws.Rows.Item(14).Delete()

# ---------------------------------------------------------------------------
# 2. Refresh the data rows (2-13) with the newly scraped opportunities.
# ---------------------------------------------------------------------------
$data = @(
    @{ Row=2;  A="1328386"; C="Sales Support Intern"; D="Panamá, Provincia de Panamá, Panamá"; E="No"; F="8 applicants";  G="6 - 18 Months"; H="Samsung Electronics Latinoamérica (Zona Libre) S.A (SELA)" },
    @{ Row=3;  A="1328376"; C="Export Manager Asistant"; D="İzmir, Türkiye"; E="No"; F="19 applicants"; G="9 - 12 Weeks";  H="BLACK LIGHT ELEKTRONIK SANAYI VE TICARET A.S." },
    @{ Row=4;  A="1328367"; C="Interior Designer"; D="Tanta, Tanta Qism 2, Tanta, Gharbia Governorate, Egypt"; E="No"; F="0 applicants";  G="9 - 12 Weeks"; H="ASG Engineering" },
    @{ Row=5;  A="1328366"; C="Architectural Engineer"; D="Tanta, Tanta Qism 2, Tanta, Gharbia Governorate, Egypt"; E="No"; F="0 applicants";  G="9 - 12 Weeks"; H="ASG Engineering" },
    @{ Row=6;  A="1328365"; C="Graphic Designer"; D="Tanta, Tanta Qism 2, Tanta, Gharbia Governorate, Egypt"; E="No"; F="0 applicants";  G="9 - 12 Weeks"; H="ASG Engineering" },
    @{ Row=7;  A="1328363"; C="Marketing Specialist"; D="Tanta, Tanta Qism 2, Tanta, Gharbia Governorate, Egypt"; E="No"; F="1 applicant";  G="9 - 12 Weeks"; H="ASG Engineering" },
    @{ Row=8;  A="1328345"; C="Design Engineer"; D="Tanta, Tanta Qism 2, Tanta, Gharbia Governorate, Egypt"; E="No"; F="0 applicants";  G="9 - 12 Weeks"; H="print shop" },
    @{ Row=9;  A="1328338"; C="Marketing Intern"; D="Sahibzada Ajit Singh Nagar, Punjab, India"; E="No"; F="4 applicants";  G="9 - 12 Weeks"; H="Crunkmart Private Limited" },
    @{ Row=10; A="1327119"; C="Marketing specialist and English instructor"; D="Heliopolis, Al Matar, El Nozha, Cairo Governorate 4470351, Egypt"; E="No"; F="10 applicants"; G="3 - 6 Months";  H="Grains nursery and preschool" },
    @{ Row=11; A="1327006"; C="[Impact Florianópolis] Global HR Innovator Internship"; D="Balneário Camboriú, SC, Brasil"; E="No"; F="44 applicants"; G="6 - 18 Months"; H="WTM do Brasil" },
    @{ Row=12; A="1326633"; C="Creative Brand & Content Coordinator"; D="Colombo, Sri Lanka"; E="No"; F="9 applicants";  G="9 - 12 Weeks"; H="Paradise Properties Ceylon" },
    @{ Row=13; A="1323494"; C="Management Trainee"; D="Polonnaruwa, Sri Lanka"; E="No"; F="3 applicants";  G="3 - 6 Months"; H="IES College International Polonnaruwa" }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = $item.A
    $ws.Cells.Item($r, 2).Value = "https://aiesec.org/opportunity/global-talent/" + $item.A
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = $item.D
    $ws.Cells.Item($r, 5).Value = $item.E
    $ws.Cells.Item($r, 6).Value = $item.F
    $ws.Cells.Item($r, 7).Value = $item.G
    $ws.Cells.Item($r, 8).Value = $item.H
}

# ---------------------------------------------------------------------------
# 3. E9 no longer carries the one-off highlight style (fillId 3 / "Yes" premium
#    call-out) - restore it to the default, unstyled cell format.
# ---------------------------------------------------------------------------
$ws.Range("E9").ClearFormats()

# ---------------------------------------------------------------------------
# 4. A few columns were re-measured (autofit-style tweak) after the refresh.
# ---------------------------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 56 - 0.83   # C: 66 -> 56
$ws.Columns.Item(4).ColumnWidth = 67 - 0.83   # D: 60 -> 67
$ws.Columns.Item(6).ColumnWidth = 16 - 0.83   # F: 17 -> 16
$ws.Columns.Item(8).ColumnWidth = 60 - 0.83   # H: 57 -> 60

"OK"
